$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record values for data rows 2-51
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 94   # AD
    $ws.Cells.Item($r, 31).Value = 68   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
